# Update "想去人数" (want-to-go count) figures in the Guangzhou Comic-Con
# info workbook, matching the regenerated gh-pages data snapshot.
#
# Sheet "展览" (Exhibitions)
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1777
$wsExhibit.Range("F11").Value = 1411
$wsExhibit.Range("F15").Value = 12729
$wsExhibit.Range("F16").Value = 12751
$wsExhibit.Range("F20").Value = 510
$wsExhibit.Range("F22").Value = 543
$wsExhibit.Range("F23").Value = 1996
$wsExhibit.Range("F24").Value = 26

# Sheet "本地生活" (Local Life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 163

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 163
$wsAll.Range("F5").Value = 1777
$wsAll.Range("F16").Value = 1411
$wsAll.Range("F21").Value = 12729
$wsAll.Range("F22").Value = 12751
$wsAll.Range("F26").Value = 510
$wsAll.Range("F28").Value = 543
$wsAll.Range("F31").Value = 1996
$wsAll.Range("F32").Value = 26
